# Auto-generated Excel COM-interop script to apply Hyperion_Profits market-data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 598.7
$ws.Range("I33").Value = 853.4
$ws.Range("J33").Value = 344
$ws.Range("K33").Value = 853.4
$ws.Range("L33").Value = 344
$ws.Range("M33").Value = -624.4
$ws.Range("N33").Value = -802

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3222.5557
$ws.Range("I40").Value = 4750.75
$ws.Range("K40").Value = 4750.75
$ws.Range("M40").Value = -4575.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 248.2
$ws.Range("I42").Value = 160.85715
$ws.Range("K42").Value = 482.57145
$ws.Range("M42").Value = -252.57145

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2426.25
$ws.Range("I100").Value = 2901.6667
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 2901.6667
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -2360.6667
$ws.Range("N100").Value = -2082

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 166667170
$ws.Range("I107").Value = 166667170
$ws.Range("K107").Value = 166667170
$ws.Range("M107").Value = -166665250

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5801.3335
$ws.Range("I113").Value = 3943.4285
$ws.Range("K113").Value = 3943.4285
$ws.Range("M113").Value = -689.4285

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9999
$ws.Range("J116").Value = 10498.833
$ws.Range("L116").Value = 10498.833
$ws.Range("N116").Value = -17382.833

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 4640.857
$ws.Range("J121").Value = 4747.6665
$ws.Range("L121").Value = 14242.9995
$ws.Range("N121").Value = -17736.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 27780540
$ws.Range("I132").Value = 27780540
$ws.Range("K132").Value = 83341620
$ws.Range("M132").Value = -83339090

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1373
$ws.Range("I135").Value = 662.4167
$ws.Range("J135").Value = 4215.3335
$ws.Range("K135").Value = 5961.7503
$ws.Range("L135").Value = 37938.0015
$ws.Range("M135").Value = -3426.7503
$ws.Range("N135").Value = -43008.0015

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 12153.23
$ws.Range("I141").Value = 14299.4
$ws.Range("K141").Value = 42898.2
$ws.Range("M141").Value = -37718.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2006.51
$ws.Range("I32").Value = 1664.1075
$ws.Range("J32").Value = 6555.5713
$ws.Range("K32").Value = 1664.1075
$ws.Range("L32").Value = 6555.5713
$ws.Range("M32").Value = -1377.1075
$ws.Range("N32").Value = -7129.5713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3110.8572
$ws.Range("I61").Value = 2641.5454
$ws.Range("K61").Value = 2641.5454
$ws.Range("M61").Value = -2429.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 706112.2
$ws.Range("I97").Value = 1081972.6
$ws.Range("K97").Value = 1081972.6
$ws.Range("M97").Value = -1081476.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 9263341
$ws.Range("I102").Value = 10420008
$ws.Range("K102").Value = 10420008
$ws.Range("M102").Value = -10418386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3969196.5
$ws.Range("I110").Value = 3969196.5
$ws.Range("K110").Value = 3969196.5
$ws.Range("M110").Value = -3967151.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2766.2163
$ws.Range("I132").Value = 2198.5356
$ws.Range("K132").Value = 6595.6068
$ws.Range("M132").Value = -4065.6068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3110.8572
$ws.Range("I136").Value = 2641.5454
$ws.Range("K136").Value = 7924.6362
$ws.Range("M136").Value = -5374.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19016.6
$ws.Range("I31").Value = 1485.3243
$ws.Range("K31").Value = 1485.3243
$ws.Range("M31").Value = -1190.3243

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 19016.6
$ws.Range("I34").Value = 1485.3243
$ws.Range("K34").Value = 1485.3243
$ws.Range("M34").Value = -1283.3243

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8864.625
$ws.Range("I58").Value = 11646.4
$ws.Range("K58").Value = 11646.4
$ws.Range("M58").Value = -11443.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 51466.668
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 51466.668
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 51466.668
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -59146.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1893.5333
$ws.Range("I122").Value = 1605.52
$ws.Range("J122").Value = 3333.6
$ws.Range("K122").Value = 4816.559999999999
$ws.Range("L122").Value = 10000.8
$ws.Range("M122").Value = -2366.559999999999
$ws.Range("N122").Value = -14900.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8864.625
$ws.Range("I136").Value = 11646.4
$ws.Range("K136").Value = 34939.2
$ws.Range("M136").Value = -32389.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 167871.33
$ws.Range("I5").Value = 1445.8
$ws.Range("J5").Value = 999999
$ws.Range("K5").Value = 4337.4
$ws.Range("L5").Value = 2999997
$ws.Range("M5").Value = -4225.4
$ws.Range("N5").Value = -3000221

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 284
$ws.Range("I8").Value = 284
$ws.Range("K8").Value = 852
$ws.Range("M8").Value = -713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 500.33334
$ws.Range("I36").Value = 500.33334
$ws.Range("K36").Value = 1501.00002
$ws.Range("M36").Value = -1332.00002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 9631.666999999999
$ws.Range("J101").Value = 9947.5
$ws.Range("L101").Value = 29842.5
$ws.Range("N101").Value = -34710.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 23810586
$ws.Range("J121").Value = 1142.6154
$ws.Range("L121").Value = 3427.8462
$ws.Range("N121").Value = -6047.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 167871.33
$ws.Range("I135").Value = 1445.8
$ws.Range("J135").Value = 999999
$ws.Range("K135").Value = 13012.2
$ws.Range("L135").Value = 8999991
$ws.Range("M135").Value = -10477.2
$ws.Range("N135").Value = -9005061

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3665.6667
$ws.Range("I138").Value = 3732.2
$ws.Range("J138").Value = 3333
$ws.Range("K138").Value = 11196.6
$ws.Range("L138").Value = 9999
$ws.Range("M138").Value = -6056.599999999999
$ws.Range("N138").Value = -20279

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1701550.6
$ws.Range("I97").Value = 2381560.8
$ws.Range("J97").Value = 1525
$ws.Range("K97").Value = 2381560.8
$ws.Range("L97").Value = 1525
$ws.Range("M97").Value = -2381064.8
$ws.Range("N97").Value = -2517

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5037180.5
$ws.Range("I102").Value = 11112829
$ws.Range("J102").Value = 1463270
$ws.Range("K102").Value = 11112829
$ws.Range("L102").Value = 1463270
$ws.Range("M102").Value = -11111207
$ws.Range("N102").Value = -1466514

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6113.8667
$ws.Range("I7").Value = 3160.1
$ws.Range("K7").Value = 3160.1
$ws.Range("M7").Value = -3048.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4967.25
$ws.Range("I46").Value = 1088.8889
$ws.Range("J46").Value = 8140.4546
$ws.Range("K46").Value = 1088.8889
$ws.Range("L46").Value = 8140.4546
$ws.Range("M46").Value = -900.8888999999999
$ws.Range("N46").Value = -8516.454600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4643.04
$ws.Range("I122").Value = 2888.4443
$ws.Range("K122").Value = 8665.332900000001
$ws.Range("M122").Value = -6215.332900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6113.8667
$ws.Range("I126").Value = 3160.1
$ws.Range("K126").Value = 9480.299999999999
$ws.Range("M126").Value = -7010.299999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 55196.562
$ws.Range("I136").Value = 85692.28999999999
$ws.Range("J136").Value = 6403.4
$ws.Range("K136").Value = 257076.87
$ws.Range("L136").Value = 19210.2
$ws.Range("M136").Value = -254526.87
$ws.Range("N136").Value = -24310.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1769
$ws.Range("I100").Value = 2447.2
$ws.Range("K100").Value = 4894.4
$ws.Range("M100").Value = -4353.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2233.8286
$ws.Range("I136").Value = 1679
$ws.Range("J136").Value = 4915.5
$ws.Range("K136").Value = 5037
$ws.Range("L136").Value = 14746.5
$ws.Range("M136").Value = -2487
$ws.Range("N136").Value = -19846.5
